# Applies the "Updated symbol list" price/volume/listing refresh described in the commit.
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h). D/E hold numeric-looking text that must
# stay literal text (exact digits/trailing zeros), so Set-TextValue forces a leading
# apostrophe (Excel's text quote-prefix) before assigning, avoiding numeric auto-conversion.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $ws.Cells.Item($row, $col).Value = "'" + $text
}

# Row 2  (BNB)
Set-TextValue 2 4 "246.94"
Set-TextValue 2 5 "0.77%"

# Row 3  (OKB)
Set-TextValue 3 4 "26.45"
Set-TextValue 3 5 "5.11%"

# Row 4  (HuobiToken)
Set-TextValue 4 4 "5.077"
Set-TextValue 4 5 "1.88%"

# Row 5  (Cronos)
Set-TextValue 5 4 "0.05599"
Set-TextValue 5 5 "-0.28%"

# Row 6  (KuCoinToken)
Set-TextValue 6 4 "6.490"
Set-TextValue 6 5 "-0.86%"

# Row 7  (MXToken)
Set-TextValue 7 4 "0.8127"
Set-TextValue 7 5 "0.13%"

# Row 8  (FTXToken)
Set-TextValue 8 4 "0.8445"
Set-TextValue 8 5 "0.54%"

# Row 9  (MandalaExchangeToken)
$ws.Cells.Item(9, 2).Value = "MandalaExchangeToken"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue 9 4 "0.06985"
Set-TextValue 9 5 "0.69%"

# Row 10  (LiechtensteinCryptoassetsExchange)
$ws.Cells.Item(10, 2).Value = "LiechtensteinCryptoassetsExchange"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue 10 4 "0.03198"
Set-TextValue 10 5 "-1.54%"

# Row 11  (BitrueCoin)
Set-TextValue 11 4 "0.02848"
Set-TextValue 11 5 "0.07%"

# Row 12  (BitMartToken)
Set-TextValue 12 4 "0.09389"
Set-TextValue 12 5 "-0.21%"

# Row 13  (BitForexToken)
Set-TextValue 13 4 "0.001515"
Set-TextValue 13 5 "0.05%"

# Row 14  (One)
$ws.Cells.Item(14, 2).Value = "One"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue 14 4 "0.0005957"
Set-TextValue 14 5 "-0.35%"

# Row 15  (TigerCash)
$ws.Cells.Item(15, 2).Value = "TigerCash"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue 15 4 "0.006149"
Set-TextValue 15 5 "-1.64%"

# Row 16  (LEO)
$ws.Cells.Item(16, 2).Value = "LEO"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue 16 4 "3.603"
Set-TextValue 16 5 "2.94%"

# Row 17  (GateToken)
$ws.Cells.Item(17, 2).Value = "GateToken"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue 17 4 "3.015"
Set-TextValue 17 5 "0.34%"

# Row 18  (BTSEToken)
$ws.Cells.Item(18, 2).Value = "BTSEToken"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue 18 4 "2.055"
Set-TextValue 18 5 "-1.73%"

# Row 19  (BitpandaEcosystemToken)
$ws.Cells.Item(19, 2).Value = "BitpandaEcosystemToken"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue 19 4 "0.3156"
Set-TextValue 19 5 "-1.26%"

# Row 20  (WazirX)
$ws.Cells.Item(20, 2).Value = "WazirX"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue 20 4 "0.1339"
Set-TextValue 20 5 "0.15%"

# Row 21  (ProBitToken)
Set-TextValue 21 5 "0.43%"

# Row 22  (MCDex)
Set-TextValue 22 4 "3.753"
Set-TextValue 22 5 "0.18%"

# Row 23  (CoinExToken)
Set-TextValue 23 4 "0.04655"
Set-TextValue 23 5 "-0.47%"

# Row 25  (BitKan)
Set-TextValue 25 4 "0.001243"
Set-TextValue 25 5 "0.14%"

# Row 26  (HotbitToken)
Set-TextValue 26 4 "0.004590"
Set-TextValue 26 5 "1.36%"

# Row 27  (NitroEx)
Set-TextValue 27 4 "0.00009599"
Set-TextValue 27 5 "-1.01%"

# Row 28  (UpBots)
Set-TextValue 28 5 "1.65%"

# Row 40  (IDEX)
Set-TextValue 40 5 "0.67%"

# Row 41  (BKEXToken)
$ws.Cells.Item(41, 2).Value = "BKEXToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue 41 4 "0.1357"
Set-TextValue 41 5 "-0.62%"

# Row 42  (CEJI)
$ws.Cells.Item(42, 2).Value = "CEJI"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue 42 4 "0.002660"
Set-TextValue 42 5 "-2.31%"

# Row 43  (KickToken)
$ws.Cells.Item(43, 2).Value = "KickToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue 43 4 "0.003411"
Set-TextValue 43 5 "-45.41%"

# Row 44  (LocalTraders)
Set-TextValue 44 4 "0.008944"
Set-TextValue 44 5 "10.86%"

# Row 45  (CoinLion)
Set-TextValue 45 4 "0.00005275"
Set-TextValue 45 5 "0.06%"

# Row 46  (Kangarootoken)
Set-TextValue 46 5 "-0.02%"

# Row 47  (CoinbaseStockToken)
Set-TextValue 47 5 "-38.90%"

# Row 48  (BOLO)
Set-TextValue 48 4 "0.002623"
Set-TextValue 48 5 "28.39%"

# Row 49  (CryptobidCoin)
Set-TextValue 49 5 "-0.02%"

# Row 50  (SpecialPowerGold)
Set-TextValue 50 5 "-0.02%"
